$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell value corrections ("J'ai corrigé l'accès au joueur") ---

# Row 6 ("Création classe BDService"): hours 3 -> 1, status "En cours" -> "Terminé"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "Terminé"

# Row 9 ("Fonctionnalités de connexion"): assigned dev Tommy Gingras -> Anthony Gauthier,
# hours 8 -> 1.5, status "Attribué" -> "Terminé"
$ws.Range("D9").Value = "Anthony Gauthier"
$ws.Range("E9").Value = 1.5
$ws.Range("F9").Value = "Terminé"

# Row 12 ("Fonctionnalité choix de la profession et du nom"): add description, status -> "En cours"
$ws.Range("C12").Value = "permet la création d'un personnage"
$ws.Range("F12").Value = "En cours"

# Row 13 ("Fonctionnalité changement de personnage"): add description, status -> "En cours"
$ws.Range("C13").Value = "permet la création d'un personnage"
$ws.Range("F13").Value = "En cours"

# Row 15 ("Fonctionnalité supprimer un personnage"): status -> "En cours"
$ws.Range("F15").Value = "En cours"

# --- View state ---
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("C17").Select()

# --- Conditional formatting: highlight "Terminé" cells/rows ---
# (added in this order so the resulting priorities come out as 1,2,3,4)

$rAll = $ws.Range("A2:F19")
$fcAll = $rAll.FormatConditions.Add(9, $null, $null, $null, "Terminé", 0)
$fcAll.Font.Color = 393372
$fcAll.Interior.Color = 13551615

$rRow2 = $ws.Range("A2:F2")
$fcRow2Equal = $rRow2.FormatConditions.Add(1, 3, '"Terminé"')
$fcRow2Equal.Font.Color = 393372
$fcRow2Equal.Interior.Color = 13551615

$rF2 = $ws.Range("F2")
$fcF2 = $rF2.FormatConditions.Add(9, $null, $null, $null, "Terminé", 0)
$fcF2.Font.Color = 393372
$fcF2.Interior.Color = 13551615

$fcRow2Contains = $rRow2.FormatConditions.Add(9, $null, $null, $null, "Terminé", 0)
$fcRow2Contains.Font.Color = 393372
$fcRow2Contains.Interior.Color = 13551615
